$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.7509363295880149
$ws1.Range("C2").Value = 0.7463235294117647
$ws1.Range("D2").Value = 0.7602996254681648
$ws1.Range("E2").Value = 0.7532467532467533
$ws1.Range("F2").Value = 0.7574626865671642
$ws1.Range("G2").Value = 0.7597524111127105
$ws1.Range("H2").Value = 0.7959187251890193
$ws1.Range("I2").Value = 406
$ws1.Range("J2").Value = 138
$ws1.Range("K2").Value = 396
$ws1.Range("L2").Value = 128

# --- Sheet "Classification Report" ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.7557251908396947
$ws2.Range("C2").Value = 0.7415730337078652
$ws2.Range("D2").Value = 0.7485822306238186

$ws2.Range("B3").Value = 0.7463235294117647
$ws2.Range("C3").Value = 0.7602996254681648
$ws2.Range("D3").Value = 0.7532467532467533

$ws2.Range("B4").Value = 0.7509363295880149
$ws2.Range("C4").Value = 0.7509363295880149
$ws2.Range("D4").Value = 0.7509363295880149
$ws2.Range("E4").Value = 0.7509363295880149

$ws2.Range("B5").Value = 0.7510243601257297
$ws2.Range("C5").Value = 0.7509363295880149
$ws2.Range("D5").Value = 0.7509144919352859

$ws2.Range("B6").Value = 0.7510243601257297
$ws2.Range("C6").Value = 0.7509363295880149
$ws2.Range("D6").Value = 0.7509144919352858

# --- Sheet "Confusion Matrix" ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 396
$ws3.Range("C2").Value = 138
$ws3.Range("B3").Value = 128
$ws3.Range("C3").Value = 406
